$d = $word.ActiveDocument

# Locate the "Science Reviewer" heading paragraph, then insert a brand new
# empty paragraph immediately after it (before whatever currently follows).
# The new paragraph keeps the same "NoSpacing" style, stays centered, and
# carries bold paragraph-mark formatting, but holds no run/text of its own.
$headingPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Science Reviewer") {
        $headingPara = $para
        break
    }
}
if ($headingPara -eq $null) {
    $headingPara = $d.Paragraphs(1)
}

$insertionPoint = $headingPara.Next().Range
$insertionPoint.Collapse(1)

$newParagraphXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="2048">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:pPr>' +
    '<w:pStyle w:val="NoSpacing"/>' +
    '<w:jc w:val="center"/>' +
    '<w:rPr><w:b/><w:bCs/></w:rPr>' +
    '</w:pPr>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$insertionPoint.InsertXML($newParagraphXml)
